$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 0.4214
$ws.Range("F10").Value = 0.851

$ws.Range("F10").Select()
